$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI LR-pair statistics (Flt3l-Flt3) per "Natmi following Dr Hou advice"
# Maps row number -> column letter -> new value
$updates = @{
    2 = @{ "E" = 2; "G" = 7.9131435; "H" = 15.826287; "I" = 0.2753713503351227; "J" = 0.2124968264014718; "K" = 3; "M" = 1.181706333333333; "N" = 3.545119; "O" = 0.2978722134968806; "P" = 0.2978722134968805; "Q" = 9.3510117905255; "R" = 56.106070743153; "S" = 0.08202547365794798; "T" = 0.06329690004126877 }
    3 = @{ "E" = 2; "G" = 7.9131435; "H" = 15.826287; "I" = 0.2753713503351227; "J" = 0.2124968264014718; "K" = 3; "M" = 2.785452333333333; "N" = 8.356357; "O" = 0.7021277865031195; "P" = 0.7021277865031195; "Q" = 22.0416840260765; "R" = 132.250104156459; "S" = 0.1933458766771748; "T" = 0.149199926360203 }
    4 = @{ "E" = 3; "G" = 3.299743666666667; "H" = 9.899231; "I" = 0.114828559503536; "J" = 0.1329152675744518; "K" = 3; "M" = 1.181706333333333; "N" = 3.545119; "O" = 0.2978722134968806; "P" = 0.2978722134968805; "Q" = 3.899327989276555; "R" = 35.093951903489; "S" = 0.03420423719197653; "T" = 0.03959176495993211 }
    5 = @{ "E" = 3; "G" = 3.299743666666667; "H" = 9.899231; "I" = 0.114828559503536; "J" = 0.1329152675744518; "K" = 3; "M" = 2.785452333333333; "N" = 8.356357; "O" = 0.7021277865031195; "P" = 0.7021277865031195; "Q" = 9.191278695718555; "R" = 82.721508261467; "S" = 0.08062432231155947; "T" = 0.09332350261451969 }
    6 = @{ "E" = 3; "G" = 3.768461; "H" = 11.305383; "I" = 0.1311395647324286; "J" = 0.1517954279960391; "K" = 3; "M" = 1.181706333333333; "N" = 3.545119; "O" = 0.2978722134968806; "P" = 0.2978722134968805; "Q" = 4.453214230619666; "R" = 40.078928075577; "S" = 0.03906283242386598; "T" = 0.04521564013588652 }
    7 = @{ "E" = 3; "G" = 3.768461; "H" = 11.305383; "I" = 0.1311395647324286; "J" = 0.1517954279960391; "K" = 3; "M" = 2.785452333333333; "N" = 8.356357; "O" = 0.7021277865031195; "P" = 0.7021277865031195; "Q" = 10.49686848552567; "R" = 94.47181636973099; "S" = 0.09207673230856267; "T" = 0.1065797878601526 }
    8 = @{ "E" = 3; "G" = 4.346190666666668; "H" = 13.038572; "I" = 0.1512441159058859; "J" = 0.1750666578210727; "K" = 3; "M" = 1.181706333333333; "N" = 3.545119; "O" = 0.2978722134968806; "P" = 0.2978722134968805; "Q" = 5.135921036674223; "R" = 46.22328933006801; "S" = 0.04505141958326499; "T" = 0.05214749287466389 }
    9 = @{ "E" = 3; "G" = 4.346190666666668; "H" = 13.038572; "I" = 0.1512441159058859; "J" = 0.1750666578210727; "K" = 3; "M" = 2.785452333333333; "N" = 8.356357; "O" = 0.7021277865031195; "P" = 0.7021277865031195; "Q" = 12.10610693357822; "R" = 108.954962402204; "S" = 0.1061926963226209; "T" = 0.1229191649464088 }
    10 = @{ "E" = 3; "G" = 5.590836; "H" = 16.772508; "I" = 0.1945568229392297; "J" = 0.22520157259838; "K" = 3; "M" = 1.181706333333333; "N" = 3.545119; "O" = 0.2978722134968806; "P" = 0.2978722134968805; "Q" = 6.606726309828001; "R" = 59.460536788452; "S" = 0.05795307149982902; "T" = 0.06708129091285787 }
    11 = @{ "E" = 3; "G" = 5.590836; "H" = 16.772508; "I" = 0.1945568229392297; "J" = 0.22520157259838; "K" = 3; "M" = 2.785452333333333; "N" = 8.356357; "O" = 0.7021277865031195; "P" = 0.7021277865031195; "Q" = 15.573007181484; "R" = 140.157064633356; "S" = 0.1366037514394007; "T" = 0.1581202816855221 }
    12 = @{ "E" = 2; "G" = 3.817888; "H" = 7.635776; "I" = 0.1328595865837971; "J" = 0.1025242476085847; "K" = 3; "M" = 1.181706333333333; "N" = 3.545119; "O" = 0.2978722134968806; "P" = 0.2978722134968805; "Q" = 4.511622429557333; "R" = 27.069734577344; "S" = 0.03957517913999609; "T" = 0.03053912457227138 }
    13 = @{ "E" = 2; "G" = 3.817888; "H" = 7.635776; "I" = 0.1328595865837971; "J" = 0.1025242476085847; "K" = 3; "M" = 2.785452333333333; "N" = 8.356357; "O" = 0.7021277865031195; "P" = 0.7021277865031195; "Q" = 10.63454503800533; "R" = 63.807270228032; "S" = 0.09328440744380098; "T" = 0.0719851230363133 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}